$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.981.65"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "3.472.48"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'587.57"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").Value = "'137.72"
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("D7").Value = "3.467.24"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -4.37%  "
$ws.Range("E10").Value = "  -6.27%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("E11").Value = "  -7.09%  "
$ws.Range("D12").Value = "'0.380"
$ws.Range("E12").Value = "  -6.87%  "
$ws.Range("D13").Value = "4.058.77"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "'0.0000181"
$ws.Range("E14").Value = "  -6.91%  "
$ws.Range("D15").Value = "'26.55"
$ws.Range("E15").Value = "  -7.48%  "
$ws.Range("D16").Value = "3.469.81"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "'0.116"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "64.910.85"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").Value = "'9.63"
$ws.Range("E19").Value = "  -10.72%  "
$ws.Range("D20").Value = "'5.75"
$ws.Range("E20").Value = "  -6.63%  "
$ws.Range("E21").Value = "  -6.13%  "
$ws.Range("D22").Value = "'388.04"
$ws.Range("E22").Value = "  -8.30%  "
$ws.Range("D23").Value = "'0.554"
$ws.Range("E23").Value = "  -5.93%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'72.35"
$ws.Range("E25").Value = "  -5.85%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.611.79"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.74"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "'0.0000108"
$ws.Range("E28").Value = "  -4.73%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'7.36"
$ws.Range("E30").Value = "  -7.10%  "
$ws.Range("D31").Value = "'8.16"
$ws.Range("E31").Value = "  -8.72%  "
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "  -10.32%  "
$ws.Range("D33").Value = "3.484.38"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -7.37%  "
$ws.Range("D36").Value = "'22.96"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("D37").Value = "'171.59"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  -10.05%  "
$ws.Range("E39").Value = "  -9.38%  "
$ws.Range("D40").Value = "'1.46"
$ws.Range("E40").Value = "  -10.35%  "
$ws.Range("E41").Value = "  -9.44%  "
$ws.Range("D42").Value = "'0.0773"
$ws.Range("E42").Value = "  -5.17%  "
$ws.Range("D43").Value = "'0.810"
$ws.Range("E43").Value = "  -5.13%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'42.24"
$ws.Range("E45").Value = "  -7.16%  "
$ws.Range("E46").Value = "  -12.88%  "
$ws.Range("D47").Value = "'23.70"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").Value = "'1.61"
$ws.Range("E48").Value = "  -9.14%  "
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").Value = "  -6.15%  "
$ws.Range("D51").Value = "2.220.15"
$ws.Range("E51").Value = "  -3.97%  "
